$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Control 30)
$ws.Range("D2").Value = 0.9999999728420853
$ws.Range("E2").Value = 0.9999999728420853

# Row 3 (Control 11)
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = 0.1659400941549308
$ws.Range("E3").Value = 0.1659400941549308

# Row 4 (Control 3)
$ws.Range("D4").Value = [double]"1.723247976499945E-24"
$ws.Range("E4").Value = [double]"1.723247976499945E-24"

# Row 5 (Control 38)
$ws.Range("C5").Value = $false
$ws.Range("D5").Value = 0.9037415003873736
$ws.Range("E5").Value = 0.9037415003873736

# Row 6 (Control 29)
$ws.Range("D6").Value = 0.8675567725311398
$ws.Range("E6").Value = 0.8675567725311398

# Row 7 (MDD 37)
$ws.Range("D7").Value = 0.9999999979084091
$ws.Range("E7").Value = [double]"2.091590900477058E-09"

# Row 8 (MDD 24)
$ws.Range("D8").Value = 0.1248569847125248
$ws.Range("E8").Value = 0.8751430152874752

# Row 9 (MDD 6)
$ws.Range("D9").Value = [double]"4.574949830204023E-06"
$ws.Range("E9").Value = 0.9999954250501698

# Row 10 (MDD 54)
$ws.Range("D10").Value = 0.001145430908042959
$ws.Range("E10").Value = 0.998854569091957

# Row 11 (MDD 21)
$ws.Range("D11").Value = 0.03389465263475766
$ws.Range("E11").Value = 0.9661053473652423
$ws.Range("F11").Value = 4.649733543395996
